$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.245.08"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "3.082.20"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D5").Value = "'547.35"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").Value = "'137.68"
$ws.Range("E6").Value = "  -6.79%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "3.083.79"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "'6.48"
$ws.Range("E11").Value = "  -3.73%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "'34.55"
$ws.Range("E13").Value = "  -5.44%  "
$ws.Range("D14").Value = "'0.0000217"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "3.568.53"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "63.247.56"
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("D17").Value = "'0.111"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "3.077.98"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'504.07"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'6.63"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'13.38"
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("D22").Value = "'0.698"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'7.13"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "'77.56"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").Value = "'12.22"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'2.71"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("D28").Value = "'8.22"
$ws.Range("E28").Value = "  -5.38%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'1.93"
$ws.Range("E30").Value = "  -8.42%  "
$ws.Range("D31").Value = "'26.23"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.50"
$ws.Range("E33").Value = "  -6.99%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'57.60"
$ws.Range("E34").Value = "  +8.82%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'518.11"
$ws.Range("E35").Value = "  -9.81%  "
$ws.Range("D36").Value = "'5.93"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "'5.15"
$ws.Range("E37").Value = "  -6.15%  "
$ws.Range("D38").Value = "'0.0404"
$ws.Range("E38").Value = "  -7.15%  "
$ws.Range("D39").Value = "3.053.54"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "'0.0784"
$ws.Range("E40").Value = "  -3.97%  "
$ws.Range("D41").Value = "'0.119"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").Value = "'2.68"
$ws.Range("E42").Value = "  -7.45%  "
$ws.Range("D43").Value = "'8.02"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").Value = "'0.252"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'121.64"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.02"
$ws.Range("E47").Value = "  -6.96%  "
$ws.Range("B48").Value = "CoreDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D48").Value = "'2.52"
$ws.Range("E48").Value = "  +74.23%  "
$ws.Range("D49").Value = "'24.05"
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("D50").Value = "'0.107"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "0.0₃0500"
$ws.Range("E51").Value = "  -5.17%  "
